$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I1").Value = "bothma-desktop"
$ws.Range("I2").Value = "D:\MS2LiveImaging\SourcePath"
$ws.Range("I3").Value = "D:\MS2LiveImaging\Code\MS2Pausing\FISHCode"
$ws.Range("I4").Value = "C:\Users\bothma\Dropbox\MS2Pausing"
$ws.Range("I6").Value = "C:\Users\bothma\Dropbox\MS2Pausing"
$ws.Range("I7").Value = "D:\MS2LiveImaging\Code\MS2Pausing\MS2Code"
$ws.Range("I8").Value = "D:\MS2LiveImaging\SchnitzcellsFolder"

$ws.Range("I1:I8").Select()
